# Add "SteadyStateTime" / "SteadyStateTimeUnit" columns to the Scenarios
# sheet (ScenarioConfiguration data), between the existing "SteadyState"
# column (G) and "ModelFile" column (which shifts from H to J).
#
# Resulting layout (row 1 = headers):
#   A Scenario_name   B IndividualId   C ModelParameterSheets
#   D ApplicationProtocol   E SimulationTime   F SimulationTimeUnit
#   G SteadyState   H SteadyStateTime (NEW)   I SteadyStateTimeUnit (NEW)
#   J ModelFile
#
# Only the TestScenario2 row (row 3) gets steady-state-time data
# (500 min); TestScenario (row 2) leaves the new cells blank, matching
# the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at H:I, pushing the old "ModelFile"
# column (and its data) from H to J.
$ws.Columns("H:I").Insert()

# New header cells.
$ws.Range("H1").Value = "SteadyStateTime"
$ws.Range("I1").Value = "SteadyStateTimeUnit"

# New data for the TestScenario2 row (row 3) only.
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = "min"

# Reflect the active selection recorded in the saved workbook.
$ws.Range("I4").Select()
